# Apply edit: insert two new rows at row 92 (pushing existing rows 92-200 down
# to 94-202), and populate the two new rows (92 and 93) with their data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows before the current row 92. This shifts the existing
# rows 92..200 down to 94..202, matching the target workbook layout.
$ws.Rows.Item(92).Insert()
$ws.Rows.Item(92).Insert()

# Populate newly inserted row 92.
$ws.Range("A92").Value() = 10
$ws.Range("B92").Value() = "Vega Modelo de Temuco"
$ws.Range("C92").Value() = "La Araucanía"
$ws.Range("D92").Value() = 44413
$ws.Range("E92").Value() = 9
$ws.Range("F92").Value() = 100112032
$ws.Range("G92").Value() = "Zapallo italiano"
$ws.Range("H92").Value() = "Sin especificar"
$ws.Range("I92").Value() = "Primera"
$ws.Range("J92").Value() = 450
$ws.Range("K92").Value() = 10000
$ws.Range("L92").Value() = 12000
$ws.Range("M92").Value() = 11111
$ws.Range("N92").Value() = "$/caja 60 unidades"
$ws.Range("O92").Value() = "Región de Arica y Parinacota"
$ws.Range("P92").Value() = 185
$ws.Range("Q92").Value() = 60
$ws.Range("R92").Value() = "Hortaliza"

# Populate newly inserted row 93.
$ws.Range("A93").Value() = 10
$ws.Range("B93").Value() = "Vega Modelo de Temuco"
$ws.Range("C93").Value() = "La Araucanía"
$ws.Range("D93").Value() = 44413
$ws.Range("E93").Value() = 9
$ws.Range("F93").Value() = 100112032
$ws.Range("G93").Value() = "Zapallo italiano"
$ws.Range("H93").Value() = "Sin especificar"
$ws.Range("I93").Value() = "Segunda"
$ws.Range("J93").Value() = 40
$ws.Range("K93").Value() = 9000
$ws.Range("L93").Value() = 9000
$ws.Range("M93").Value() = 9000
$ws.Range("N93").Value() = "$/caja 80 unidades"
$ws.Range("O93").Value() = "Región de Arica y Parinacota"
$ws.Range("P93").Value() = 112
$ws.Range("Q93").Value() = 80
$ws.Range("R93").Value() = "Hortaliza"

# Make sure the date cells carry the same date number format as the rest of
# column D (style should already be inherited from the Insert, but set it
# explicitly to be safe).
$ws.Range("D92").NumberFormat = $ws.Range("D91").NumberFormat
$ws.Range("D93").NumberFormat = $ws.Range("D91").NumberFormat
